$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.00210731725890954"
$ws.Range("E2").Value = [double]"0.00210731725890954"

$ws.Range("D3").Value = [double]"3.052747142492149E-28"
$ws.Range("E3").Value = [double]"3.052747142492149E-28"

$ws.Range("D4").Value = [double]"1.295550330632437E-11"
$ws.Range("E4").Value = [double]"1.295550330632437E-11"

$ws.Range("D5").Value = [double]"0.9996557658806731"
$ws.Range("E5").Value = [double]"0.9996557658806731"

$ws.Range("D6").Value = [double]"0.9999999999999631"
$ws.Range("E6").Value = [double]"0.9999999999999631"

$ws.Range("D7").Value = [double]"2.701160201811288E-06"
$ws.Range("E7").Value = [double]"0.9999972988397982"

$ws.Range("D8").Value = [double]"0.9999996923086071"
$ws.Range("E8").Value = [double]"3.076913929200842E-07"

$ws.Range("D11").Value = [double]"0.999999999999855"
$ws.Range("E11").Value = [double]"1.449951270160454E-13"
$ws.Range("F11").Value = [double]"5.172896385192871"
